# Weekly update: insert the latest week's price records for
# "Agrícola del Norte S.A. de Arica" / Zapallo italiano, Huracán
# right before the existing 2021-01-29 entries (row 436), shifting
# every subsequent record down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 436/437 (everything from old row 436 onward
# moves down by two rows; formatting of the row above - including the
# date number format on column D - is carried down automatically).
$ws.Range("A436:A437").EntireRow.Insert()

# New row 436: "Primera" quality record for the new week.
$ws.Cells.Item(436, 1).Value = 1
$ws.Cells.Item(436, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(436, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(436, 4).Value = 44890
$ws.Cells.Item(436, 5).Value = 15
$ws.Cells.Item(436, 6).Value = 100112032
$ws.Cells.Item(436, 7).Value = "Zapallo italiano"
$ws.Cells.Item(436, 8).Value = "Huracán"
$ws.Cells.Item(436, 9).Value = "Primera"
$ws.Cells.Item(436, 10).Value = 490
$ws.Cells.Item(436, 11).Value = 2500
$ws.Cells.Item(436, 12).Value = 3000
$ws.Cells.Item(436, 13).Value = 2694
$ws.Cells.Item(436, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(436, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(436, 16).Value = 38
$ws.Cells.Item(436, 17).Value = 70
$ws.Cells.Item(436, 18).Value = "Hortaliza"

# New row 437: "Segunda" quality record for the new week.
$ws.Cells.Item(437, 1).Value = 1
$ws.Cells.Item(437, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(437, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(437, 4).Value = 44890
$ws.Cells.Item(437, 5).Value = 15
$ws.Cells.Item(437, 6).Value = 100112032
$ws.Cells.Item(437, 7).Value = "Zapallo italiano"
$ws.Cells.Item(437, 8).Value = "Huracán"
$ws.Cells.Item(437, 9).Value = "Segunda"
$ws.Cells.Item(437, 10).Value = 410
$ws.Cells.Item(437, 11).Value = 2000
$ws.Cells.Item(437, 12).Value = 2500
$ws.Cells.Item(437, 13).Value = 2220
$ws.Cells.Item(437, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(437, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(437, 16).Value = 22
$ws.Cells.Item(437, 17).Value = 100
$ws.Cells.Item(437, 18).Value = "Hortaliza"
